$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Insert two new rows before old row 6 (My Profile) to hold BS View Payments and Payer View Payments
$ws.Range("A6:C7").Insert()

$ws.Range("A6").Value = "UPA_Regression"
$ws.Range("A7").Value = "UPA_Regression"

# Write the label column first (matches shared-string insertion order)
$ws.Range("B5").Value = "Provider View Payments"
$ws.Range("B6").Value = "BS View Payments"
$ws.Range("B7").Value = "Payer View Payments"

# Then the classname column
$ws.Range("C5").Value = "test.java.TestProviderViewPayments"
$ws.Range("C6").Value = "test.java.TestBSViewPayments"
$ws.Range("C7").Value = "test.java.TestPayerViewPayments"

# Row 8 now holds My Profile (previously row 6), already correct values - no change needed
# Row 9 now holds CSR_Regression / Manage Users / TestCSRManageUsers (previously row 7), no change needed

# Add new row 10: CSR_Regression / View Payments / TestCSRViewPayments
$ws.Range("A10").Value = "CSR_Regression"
$ws.Range("B10").Value = "View Payments"
$ws.Range("C10").Value = "TestCSRViewPayments"

$ws.Range("C10").Select()
